$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B12").Value = "Tomek"
$ws.Range("C12").Value = "Done"
$ws.Range("B13").Value = "Tomek"
$ws.Range("C13").Value = "Done"

$ws.Range("B13:C13").Select()
